$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Row 7 ---
$ws.Range("A7").Value = 131263749
$ws.Range("B7").Value = 58256
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 103015
$ws.Range("F7").Value = "Kungsfågel"
$ws.Range("G7").Value = "Regulus regulus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"

# Text that looks numeric must stay plain text (not be auto-converted to a number).
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1"
$ws.Range("I7").Style = "Normal"

$ws.Range("P7").Value = "Hjärtaboda, Hjärtaboda, Sk"
$ws.Range("Q7").Value = 448813
$ws.Range("R7").Value = 6230068
$ws.Range("S7").Value = 30
$ws.Range("T7").Value = "Skåne"
$ws.Range("U7").Value = "Östra Göinge"
$ws.Range("V7").Value = "Skåne"
$ws.Range("W7").Value = "Hjärsås"

# Date-like text must stay plain text (not be auto-converted to a date serial).
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2026-02-22"
$ws.Range("Y7").Style = "Normal"

$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2026-02-22"
$ws.Range("AA7").Style = "Normal"

$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = "Martin Kornhall"
$ws.Range("AX7").Value = "Martin Kornhall"

# --- Row 8 ---
$ws.Range("A8").Value = 131263369
$ws.Range("B8").Value = 92466
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 6031
$ws.Range("F8").Value = "Blomkålssvamp"
$ws.Range("G8").Value = "Sparassis crispa"
$ws.Range("H8").Value = "(Wulfen:Fr.) Fr."

$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "6"
$ws.Range("I8").Style = "Normal"

$ws.Range("P8").Value = "Hjärtaboda, Hjärtaboda, Sk"
$ws.Range("Q8").Value = 448731
$ws.Range("R8").Value = 6230010
$ws.Range("S8").Value = 40
$ws.Range("T8").Value = "Skåne"
$ws.Range("U8").Value = "Östra Göinge"
$ws.Range("V8").Value = "Skåne"
$ws.Range("W8").Value = "Hjärsås"

$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2020-09-26"
$ws.Range("Y8").Style = "Normal"

$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2020-09-26"
$ws.Range("AA8").Style = "Normal"

$ws.Range("AC8").Value = "Växte vid foten av gamla grova tallar"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = "Martin Kornhall"
$ws.Range("AX8").Value = "Martin Kornhall"
